$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range('D2').Value = '43.881.54'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '2.317.00'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''232.32'
$ws.Range('E5').Value = '  +0.94%  '
$ws.Range('E6').Value = '  +1.62%  '
$ws.Range('D7').Value = '''65.35'
$ws.Range('E7').Value = '  +3.25%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.444'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').Value = '''0.0970'
$ws.Range('E10').Value = '  -4.57%  '
$ws.Range('D11').Value = '''56.97'
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '''26.65'
$ws.Range('E12').Value = '  +2.80%  '
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').Value = '2.652.49'
$ws.Range('E14').Value = '  +2.35%  '
$ws.Range('D15').Value = '''15.33'
$ws.Range('E15').Value = '  -1.68%  '
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '''0.836'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '2.317.05'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('D19').Value = '43.786.47'
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('D20').Value = '0.0₃0978'
$ws.Range('E20').Value = '  -2.47%  '
$ws.Range('D21').Value = '''73.65'
$ws.Range('E21').Value = '  +0.68%  '
$ws.Range('D22').Value = '''6.18'
$ws.Range('E22').Value = '  +2.67%  '
$ws.Range('D23').Value = '''249.57'
$ws.Range('E23').Value = '  -0.98%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').Value = '''3.73'
$ws.Range('E25').Value = '  +12.03%  '
$ws.Range('E26').Value = '  +0.50%  '
$ws.Range('E27').Value = '  -2.64%  '
$ws.Range('D28').Value = '''9.83'
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '''22.19'
$ws.Range('E29').Value = '  +7.24%  '
$ws.Range('B30').Value = 'Monero'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D30').Value = '''174.50'
$ws.Range('E30').Value = '  +1.66%  '
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.44'
$ws.Range('E31').Value = '  +4.85%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '''0.132'
$ws.Range('E32').Value = '  -4.14%  '
$ws.Range('E33').Value = '  +1.51%  '
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('E35').Value = '  -0.09%  '
$ws.Range('D36').Value = '''4.97'
$ws.Range('E36').Value = '  +2.78%  '
$ws.Range('D37').Value = '''6.53'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '''3.65'
$ws.Range('E38').Value = '  -3.81%  '
$ws.Range('D39').Value = '''2.39'
$ws.Range('E39').Value = '  +4.20%  '
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').Value = '''9.04'
$ws.Range('E41').Value = '  +9.76%  '
$ws.Range('B42').Value = 'BinanceUSD'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D42').Value = '''1.00'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '''17.39'
$ws.Range('E43').Value = '  -0.04%  '
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '''98.52'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = '''4.45'
$ws.Range('E45').Value = '  +1.31%  '
$ws.Range('B46').Value = 'Celestia'
$ws.Range('C46').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').Value = '''10.34'
$ws.Range('E46').Value = '  +7.53%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '''0.0954'
$ws.Range('E47').Value = '  -1.65%  '
$ws.Range('E48').Value = '  -0.49%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '''1.12'
$ws.Range('E49').Value = '  +4.48%  '
$ws.Range('B50').Value = 'TerraClassic'
$ws.Range('C50').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D50').Value = '''0.000208'
$ws.Range('E50').Value = '  -4.51%  '
$ws.Range('D51').Value = '1.443.34'
$ws.Range('E51').Value = '  +0.29%  '
